$wb = $excel.ActiveWorkbook

# 1. Rename sheet_2 -> シート<2> (exercises XML escaping of sheet name: < and >)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "シート<2>"

# 2. Add a new row on sheet 2 with strings that need XML escaping
#    (<, >, ", ', & and an embedded newline)
$ws1 = $wb.Worksheets.Item(1)
$ws2.Range("A2").Value = "<>`"'&"
$ws2.Range("B2").Value = "a`nb"

# 3. Change sheet 1 B3 from a date+time value to a plain date value
#    (keeps existing date number format / style)
$ws1.Range("B3").Value = "10/10/2012"
